$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 391; this shifts rows 391:423 down to 392:424
$ws.Rows.Item(391).Insert()

# Populate the newly inserted row 391 with the new data record
$ws.Cells.Item(391, 1).Value = 10
$ws.Cells.Item(391, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(391, 3).Value = "La Araucanía"
$ws.Cells.Item(391, 4).Value = 45013
$ws.Cells.Item(391, 5).Value = 9
$ws.Cells.Item(391, 6).Value = 100112001
$ws.Cells.Item(391, 7).Value = "Berenjena"
$ws.Cells.Item(391, 8).Value = "Sin especificar"
$ws.Cells.Item(391, 9).Value = "Primera"
$ws.Cells.Item(391, 10).Value = 50
$ws.Cells.Item(391, 11).Value = 14000
$ws.Cells.Item(391, 12).Value = 14000
$ws.Cells.Item(391, 13).Value = 14000
$ws.Cells.Item(391, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(391, 15).Value = "Región del Maule"
$ws.Cells.Item(391, 16).Value = 350
$ws.Cells.Item(391, 17).Value = 40
$ws.Cells.Item(391, 18).Value = "Hortaliza"
